$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28 (shifts existing rows 28..116 down to 29..117)
$ws.Rows.Item(28).Insert()

# Fill in the new row 28 with the new daily sales record (May 27, 2025)
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = 24156.31
$ws.Cells.Item(28, 3).Value = 5
$ws.Cells.Item(28, 4).Value = 2025
$ws.Cells.Item(28, 5).Value = "05/2025"
